# Optuna Attempt (go back with original)
$wb = $excel.ActiveWorkbook

# Sheet 1: Forecast Comparison
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$ws1.Range("L2").Value = 0.89
$ws1.Range("L3").Value = 1.01
$ws1.Range("L4").Value = 0.83
$ws1.Range("L5").Value = 1.01
$ws1.Range("L6").Value = 0.9399999999999999
$ws1.Range("L7").Value = 0.95
$ws1.Range("L8").Value = 0.89

$ws1.Range("D9").Value = 48
$ws1.Range("L9").Value = 0.89

$ws1.Range("D10").Value = 45
$ws1.Range("L10").Value = 1.19

$ws1.Range("D11").Value = 43
$ws1.Range("L11").Value = 1.19

$ws1.Range("D12").Value = 41
$ws1.Range("L12").Value = 0.9399999999999999

$ws1.Range("D13").Value = 41
$ws1.Range("L13").Value = 0.9399999999999999

$ws1.Range("D14").Value = 38
$ws1.Range("L14").Value = 1.17

$ws1.Range("D15").Value = 38
$ws1.Range("L15").Value = 0.91

$ws1.Range("D16").Value = 36
$ws1.Range("L16").Value = 0.8100000000000001

$ws1.Range("D17").Value = 33
$ws1.Range("L17").Value = 1.18

# Sheet 2: Summary
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "708"

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "393"

$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = "33"
